$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConDA")
$ws.Activate()

# Insert 3 new rows before row 30, using rows 24:26 as the formatting template
# (merged B/C columns, single D, merged E holding the seed, F/G/H metrics, blank merged I)
$ws.Rows("24:26").Copy()
$ws.Rows("30:32").Insert()

# Re-merge the newly inserted block's B/C/E columns (Insert() does not recreate merges)
$ws.Range("B30:B32").Merge()
$ws.Range("C30:C32").Merge()
$ws.Range("E30:E32").Merge()

# Fill in the content for the new block (rows 30-32)
$ws.Range("B30").Value = "ConDA(BLIP-2) w/ test_time_adaptation()  before validate() w/o triplet loss, using z as the input to the classifier instead of h"
$ws.Range("C30").Value = "toy training set"

$ws.Range("D30").Value = "Source: Covid, Climate" + [char]10 + "Target: Military"
$ws.Range("E30").Value = 1001
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "NA"
$ws.Range("H30").Value = 0.796

$ws.Range("D31").Value = "Source: Covid, Military" + [char]10 + "Target: Climate"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = 0.801
$ws.Range("H31").Value = "NA"

$ws.Range("D32").Value = "Source: Climate, Military" + [char]10 + "Target: Covid"
$ws.Range("F32").Value = 0.801
$ws.Range("G32").Value = "NA"
$ws.Range("H32").Value = "NA"

# Make sure I30:I32 has no value (blank, matching the template source)
$ws.Range("I30").ClearContents()

# Update the view: scroll down and select the merged I cell of the row that used to be selected (shifted by 3)
$excel.Goto($ws.Range("A25"), $true)
$ws.Range("I33:I35").Select()

$wb.Save()
